# Updated cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.622.38"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.596.03"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.53"
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("E6").Value = "  +1.43%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0617"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.43"
$ws.Range("E10").Value = "  -0.74%  "
$ws.Range("E11").Value = "  +0.55%  "
$ws.Range("D12").Value = "1.819.02"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "1.603.63"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.59"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "26.594.58"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.80"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("E21").Value = "  +3.65%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.88"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.12"
$ws.Range("E25").Value = "  -1.02%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  +0.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.24"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0505"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +0.36%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.651"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").Value = "1.281.34"
$ws.Range("E35").Value = "  -1.72%  "
$ws.Range("E36").Value = "  +0.93%  "
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("E42").Value = "  +1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.11"
$ws.Range("E44").Value = "  +2.58%  "
$ws.Range("D45").Value = "1.731.10"
$ws.Range("E45").Value = "  +0.45%  "
$ws.Range("E46").Value = "  +8.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.64"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("E50").Value = "  +4.49%  "
$ws.Range("E51").Value = "  +0.34%  "
